# Update LR-pairs data with new TPM-derived values (Amh-Amhr2, YoungD7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, followed by the 20 column values (A:T)
$rows = @(
    ,@(2, "FAPs", "Amh", "Amhr2", "ECs", 1, 0.3333333333333333, 0.294587, 0.883761, 0.2010206147344287, 0.2010206147344287, 1, 0.3333333333333333, 0.1785033333333333, 0.53551, 0.132093478502712, 0.132093478502712, 0.05258476145666667, 0.47326285311, 0.02655351225102421, 0.0265535122510242)
    ,@(3, "FAPs", "Amh", "Amhr2", "FAPs", 1, 0.3333333333333333, 0.294587, 0.883761, 0.2010206147344287, 0.2010206147344287, 2, 0.6666666666666666, 0.6758093333333334, 2.027428, 0.5001027374536356, 0.5001027374536355, 0.1990846440786667, 1.791761796708, 0.1005309597133004, 0.1005309597133004)
    ,@(4, "FAPs", "Amh", "Amhr2", "MuSCs", 1, 0.3333333333333333, 0.294587, 0.883761, 0.2010206147344287, 0.2010206147344287, 2, 0.6666666666666666, 0.3973076666666667, 1.191923, 0.2940099254493623, 0.2940099254493623, 0.1170416736003333, 1.053375062403, 0.05910205595185435, 0.05910205595185435)
    ,@(5, "FAPs", "Amh", "Amhr2", "Resolving-Mac", 1, 0.3333333333333333, 0.294587, 0.883761, 0.2010206147344287, 0.2010206147344287, 1, 0.3333333333333333, 0.09972066666666667, 0.299162, 0.07379385859429018, 0.07379385859429016, 0.02937641203133333, 0.264387708282, 0.01483408681824971, 0.01483408681824971)
    ,@(6, "MuSCs", "Amh", "Amhr2", "ECs", 3, 1, 1.032183333333333, 3.09655, 0.7043424461544411, 0.7043424461544411, 1, 0.3333333333333333, 0.1785033333333333, 0.53551, 0.132093478502712, 0.132093478502712, 0.1842481656111111, 1.6582334905, 0.09303904376964926, 0.09303904376964925)
    ,@(7, "MuSCs", "Amh", "Amhr2", "FAPs", 3, 1, 1.032183333333333, 3.09655, 0.7043424461544411, 0.7043424461544411, 2, 0.6666666666666666, 0.6758093333333334, 2.027428, 0.5001027374536356, 0.5001027374536355, 0.6975591303777777, 6.2780321734, 0.3522435854266259, 0.3522435854266259)
    ,@(8, "MuSCs", "Amh", "Amhr2", "MuSCs", 3, 1, 1.032183333333333, 3.09655, 0.7043424461544411, 0.7043424461544411, 2, 0.6666666666666666, 0.3973076666666667, 1.191923, 0.2940099254493623, 0.2940099254493623, 0.4100943517388889, 3.69084916565, 0.2070836700846887, 0.2070836700846887)
    ,@(9, "MuSCs", "Amh", "Amhr2", "Resolving-Mac", 3, 1, 1.032183333333333, 3.09655, 0.7043424461544411, 0.7043424461544411, 1, 0.3333333333333333, 0.09972066666666667, 0.299162, 0.07379385859429018, 0.07379385859429016, 0.1029300101222222, 0.9263700910999999, 0.05197614687347727, 0.05197614687347726)
    ,@(10, "Resolving-Mac", "Amh", "Amhr2", "ECs", 1, 0.3333333333333333, 0.1386863333333333, 0.416059, 0.09463693911113033, 0.09463693911113034, 1, 0.3333333333333333, 0.1785033333333333, 0.53551, 0.132093478502712, 0.132093478502712, 0.02475597278777778, 0.22280375509, 0.01250092248203856, 0.01250092248203856)
    ,@(11, "Resolving-Mac", "Amh", "Amhr2", "FAPs", 1, 0.3333333333333333, 0.1386863333333333, 0.416059, 0.09463693911113033, 0.09463693911113034, 2, 0.6666666666666666, 0.6758093333333334, 2.027428, 0.5001027374536356, 0.5001027374536355, 0.09372551847244445, 0.843529666252, 0.04732819231370931, 0.04732819231370931)
    ,@(12, "Resolving-Mac", "Amh", "Amhr2", "MuSCs", 1, 0.3333333333333333, 0.1386863333333333, 0.416059, 0.09463693911113033, 0.09463693911113034, 2, 0.6666666666666666, 0.3973076666666667, 1.191923, 0.2940099254493623, 0.2940099254493623, 0.05510114349522222, 0.4959102914570001, 0.02782419941281927, 0.02782419941281927)
    ,@(13, "Resolving-Mac", "Amh", "Amhr2", "Resolving-Mac", 1, 0.3333333333333333, 0.1386863333333333, 0.416059, 0.09463693911113033, 0.09463693911113034, 1, 0.3333333333333333, 0.09972066666666667, 0.299162, 0.07379385859429018, 0.07379385859429016, 0.01382989361755555, 0.124469042558, 0.006983624902563201, 0.006983624902563201)
)

foreach ($entry in $rows) {
    $rownum = $entry[0]
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0, $i] = $entry[$i + 1]
    }
    $ws.Range("A" + $rownum + ":T" + $rownum).Value2 = $arr
}
